$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header cell's format (H1) onto the two new header cells (I1:J1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8
$ws.Range("I3").Value = 8
$ws.Range("J3").Value = 8
$ws.Range("I4").Value = 9
$ws.Range("J4").Value = 9
$ws.Range("I5").Value = 9
$ws.Range("J5").Value = 9
